# Add a header row to the Europe temperature-change data:
#   A: "Year"   B: "Temperature change"
# Existing data (previously rows 1-49) shifts down to rows 2-50.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 1, pushing all data down.
$ws.Rows.Item(1).Insert() | Out-Null

# Fill in the new header row.
$ws.Range("A1").Value = "Year"
$ws.Range("B1").Value = "Temperature change"

# Move/collapse the selection to A2, matching the post-edit workbook state.
$ws.Range("A2").Select() | Out-Null
